$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (Hora): all data rows 2-51 change from 21 to 22
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "22"

# Column D (Price) updates
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("D2").Value = "316.28"
$ws.Range("D3").Value = "37.97"
$ws.Range("D4").Value = "5.170"
$ws.Range("D5").Value = "0.07972"
$ws.Range("D6").Value = "8.452"
$ws.Range("D7").Value = "1.932"
$ws.Range("D9").Value = "0.9426"
$ws.Range("D10").Value = "0.1248"
$ws.Range("D11").Value = "0.1936"
$ws.Range("D12").Value = "0.09021"
$ws.Range("D13").Value = "0.03425"
$ws.Range("D14").Value = "0.09531"
$ws.Range("D15").Value = "0.001368"
$ws.Range("D16").Value = "0.006049"
$ws.Range("D17").Value = "3.426"
$ws.Range("D18").Value = "4.472"
$ws.Range("D19").Value = "0.3515"
$ws.Range("D20").Value = "6.518"
$ws.Range("D21").Value = "0.1306"
$ws.Range("D22").Value = "0.2306"
$ws.Range("D23").Value = "0.04360"
$ws.Range("D24").Value = "0.001225"
$ws.Range("D26").Value = "0.0001325"
$ws.Range("D39").Value = "0.02397"
$ws.Range("D40").Value = "0.05178"
$ws.Range("D41").Value = "0.007446"
$ws.Range("D43").Value = "0.008404"
$ws.Range("D44").Value = "0.002092"
$ws.Range("D45").Value = "0.008727"
$ws.Range("D46").Value = "0.00006475"
$ws.Range("D47").Value = "0.00000000747"
$ws.Range("D48").Value = "0.002856"
$ws.Range("D50").Value = "0.00002092"
$ws.Range("D51").Value = "0.0001993"

# Column E (Volume(1h)) updates
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("E2").Value = "1.45%"
$ws.Range("E3").Value = "1.58%"
$ws.Range("E4").Value = "0.89%"
$ws.Range("E5").Value = "1.85%"
$ws.Range("E6").Value = "2.19%"
$ws.Range("E7").Value = "1.61%"
$ws.Range("E8").Value = "5.07%"
$ws.Range("E9").Value = "2.54%"
$ws.Range("E10").Value = "4.54%"
$ws.Range("E11").Value = "0.96%"
$ws.Range("E12").Value = "0.48%"
$ws.Range("E13").Value = "2.12%"
$ws.Range("E14").Value = "-0.67%"
$ws.Range("E15").Value = "-1.71%"
$ws.Range("E16").Value = "6.22%"
$ws.Range("E17").Value = "-3.06%"
$ws.Range("E18").Value = "1.19%"
$ws.Range("E19").Value = "2.15%"
$ws.Range("E20").Value = "23.51%"
$ws.Range("E21").Value = "1.82%"
$ws.Range("E22").Value = "-11.19%"
$ws.Range("E23").Value = "-0.08%"
$ws.Range("E24").Value = "-2.07%"
$ws.Range("E25").Value = "-5.65%"
$ws.Range("E26").Value = "-2.79%"
$ws.Range("E27").Value = "-0.49%"
$ws.Range("E39").Value = "4.91%"
$ws.Range("E40").Value = "2.50%"
$ws.Range("E41").Value = "-0.27%"
$ws.Range("E42").Value = "3.48%"
$ws.Range("E43").Value = "-7.34%"
$ws.Range("E44").Value = "7.06%"
$ws.Range("E45").Value = "-6.36%"
$ws.Range("E46").Value = "-2.27%"
$ws.Range("E47").Value = "-0.49%"
$ws.Range("E48").Value = "-12.87%"
$ws.Range("E49").Value = "68.11%"
$ws.Range("E50").Value = "-0.49%"
$ws.Range("E51").Value = "-0.49%"
